$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (TPM-updated values) for rows 2-7, columns A:T.
# Columns A-D are text labels (resolved through shared strings automatically
# when assigned as .Value strings), columns E-T are numeric.
$data = @(
    @("ECs",  "Tnfsf10", "Tnfrsf11b", "FAPs",  3, 1,               17.31937166666667,  51.958115,  0.9799257492601914,   0.9799257492601915,   3, 1, 1.183046666666667, 3.54914,  0.6222589862820888,  0.6222589862820888,  20.48962491901111,  184.4066242711,    0.6097676033663629,   0.609767603366363),
    @("ECs",  "Tnfsf10", "Tnfrsf11b", "MuSCs", 3, 1,               17.31937166666667,  51.958115,  0.9799257492601914,   0.9799257492601915,   3, 1, 0.718166,          2.154498, 0.3777410137179113,  0.3777410137179112,  12.43818387236333,  111.94365485127,   0.3701581458938284,   0.3701581458938284),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "FAPs",  1, 0.3333333333333333, 0.170846,        0.512538,   0.009666424266436919, 0.009666424266436919, 3, 1, 1.183046666666667, 3.54914,  0.6222589862820888,  0.6222589862820888,  0.2021187908133333, 1.81906911732,     0.006015019365005621, 0.006015019365005621),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "MuSCs", 1, 0.3333333333333333, 0.170846,        0.512538,   0.009666424266436919, 0.009666424266436919, 3, 1, 0.718166,          2.154498, 0.3777410137179113,  0.3777410137179112,  0.122695788436,     1.104262095924,    0.003651404901431299, 0.003651404901431299),
    @("MuSCs","Tnfsf10", "Tnfrsf11b", "FAPs",  2, 0.6666666666666666, 0.1839496666666667, 0.551849, 0.01040782647337163,  0.01040782647337163,  3, 1, 1.183046666666667, 3.54914,  0.6222589862820888,  0.6222589862820888,  0.2176210399844444, 1.95858935986,     0.006476363550720115, 0.006476363550720117),
    @("MuSCs","Tnfsf10", "Tnfrsf11b", "MuSCs", 2, 0.6666666666666666, 0.1839496666666667, 0.551849, 0.01040782647337163,  0.01040782647337163,  3, 1, 0.718166,          2.154498, 0.3777410137179113,  0.3777410137179112,  0.1321063963113333, 1.188957566802,    0.003931462922651512, 0.003931462922651512)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}

# Remove the now-unused former rows 8-10 (old data had a 3x3 cross of
# clusters; new data only has a 3x2 cross, so the trailing rows go away).
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null
